# New weekly price record for "Haba" (Comercializadora del Agro de Limarí)
# is inserted as row 64, pushing the existing rows 64:96 down to 65:97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64 (shifts rows 64-96 down to 65-97,
# dimension grows from A1:R96 to A1:R97 automatically).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new weekly record.
$ws.Cells.Item(64, 1).Value = 2
$ws.Cells.Item(64, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 45120
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = 100112026
$ws.Cells.Item(64, 7).Value = "Haba"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 1600
$ws.Cells.Item(64, 11).Value = 10000
$ws.Cells.Item(64, 12).Value = 12000
$ws.Cells.Item(64, 13).Value = 11000
$ws.Cells.Item(64, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(64, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 16).Value = 440
$ws.Cells.Item(64, 17).Value = 25
$ws.Cells.Item(64, 18).Value = "Hortaliza"
